# "Generate Report for Handoff" — mark the 86ae12c6-... and cdbd9a21-... files
# as handed off (status -> "Ready for handoff"), refresh their handoff
# timestamps, and record a "handback file not latest" error detail for the
# zh-cn / de-de localization rows. Mirrors the report generator's per-run
# update of the localization-status workbook.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: rows for 86ae12c6-... (row 4) and cdbd9a21-... (row 5)
# Columns: E=zh-cn status, F=de-de status, G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $readyForHandoff
$wsOverview.Range("F4").Value = $readyForHandoff
$wsOverview.Range("G4").Value = "2016-09-05 06:28:59"

$wsOverview.Range("E5").Value = $readyForHandoff
$wsOverview.Range("F5").Value = $readyForHandoff
$wsOverview.Range("G5").Value = "2016-09-05 06:28:59"

# ---------------------------------------------------------------------
# zh-cn sheet: rows for 86ae12c6-... (row 4) and cdbd9a21-... (row 5)
# Columns: C=Status, H=Latest Handoff Datetime, P=Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

$wsZhCn.Range("C4").Value = $readyForHandoff
$wsZhCn.Range("H4").Value = "2016-09-05 06:28:54"
$wsZhCn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3f2ed53abe0bbd2dfda50bf72be4e6ad6ba5d92/e2e/86ae12c6-bbac-4a52-8c0b-d6184a662e28.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da40aca201bfa52b6eb6ae394981eaeb67cffe09/e2e/86ae12c6-bbac-4a52-8c0b-d6184a662e28.md."

$wsZhCn.Range("C5").Value = $readyForHandoff
$wsZhCn.Range("H5").Value = "2016-09-05 06:28:54"
$wsZhCn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3f2ed53abe0bbd2dfda50bf72be4e6ad6ba5d92/e2e/cdbd9a21-d465-4594-b4a9-547ef1279cf4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da40aca201bfa52b6eb6ae394981eaeb67cffe09/e2e/cdbd9a21-d465-4594-b4a9-547ef1279cf4.md."

# ---------------------------------------------------------------------
# de-de sheet: rows for 86ae12c6-... (row 4) and cdbd9a21-... (row 5)
# Columns: C=Status, H=Latest Handoff Datetime, P=Error Detail
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Columns.Item(16).ColumnWidth = 39.17

$wsDeDe.Range("C4").Value = $readyForHandoff
$wsDeDe.Range("H4").Value = "2016-09-05 06:28:59"
$wsDeDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3f2ed53abe0bbd2dfda50bf72be4e6ad6ba5d92/e2e/86ae12c6-bbac-4a52-8c0b-d6184a662e28.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da40aca201bfa52b6eb6ae394981eaeb67cffe09/e2e/86ae12c6-bbac-4a52-8c0b-d6184a662e28.md."

$wsDeDe.Range("C5").Value = $readyForHandoff
$wsDeDe.Range("H5").Value = "2016-09-05 06:28:59"
$wsDeDe.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3f2ed53abe0bbd2dfda50bf72be4e6ad6ba5d92/e2e/cdbd9a21-d465-4594-b4a9-547ef1279cf4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da40aca201bfa52b6eb6ae394981eaeb67cffe09/e2e/cdbd9a21-d465-4594-b4a9-547ef1279cf4.md."
